# Append the new call-log entries to the "Sheet" worksheet (column A),
# continuing the existing log starting at row 275.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")

$entries = @(
    "15 15:03>>> 4FD889D140   Freddy Velez",
    "15 15:03>>> 4FD889D140   Freddy Velez",
    "18 14:33>>> 4FD889D140   Freddy Velez",
    "18 14:35>>> 4FD889D140   Freddy Velez",
    "18 14:35>>> 4FD889D140   Freddy Velez",
    "18 14:35>>> 4FD889D140   Freddy Velez",
    "18 14:37>>> 4FD889D140   Freddy Velez",
    "18 14:37>>> 4FD889D140   Freddy Velez",
    "18 14:38>>> 4FD889D140   Freddy Velez",
    "18 14:41>>> 4FD889D140   Freddy Velez",
    "18 14:41>>> 4FD889D140   Freddy Velez"
)

$startRow = 275
for ($i = 0; $i -lt $entries.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $entries[$i]
}
